$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")

$ws.Range("C58").Value = 6940
$ws.Range("C59").Value = 2210
$ws.Range("C60").Value = 2210

$ws.Range("C60").Select()
